$wb = $excel.ActiveWorkbook

# --- Sheet "Score Cards": append a new row of scorecard data ---
$ws1 = $wb.Worksheets.Item("Score Cards")

$newRow = @("Ocean View", 45633, 4, 5, 7, 4, 6, 4, 7, 4, 6, 5, 5, 6, 4, 5, 3, 6, 4, 6, 0, 34, 7)
for ($i = 0; $i -lt $newRow.Length; $i++) {
    $ws1.Cells.Item(4, $i + 1).Value = $newRow[$i]
}

# Copy the formatting from the row above so the new row picks up the same
# cell styles (date number format on column B, etc.) rather than defaults.
$ws1.Range("A3:W3").Copy()
$ws1.Range("A4:W4").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Grow Table1 so the new row is included in the table range.
$lo1 = $ws1.ListObjects.Item(1)
$lo1.Resize($ws1.Range("A1:W4"))

# Update the selection shown on this sheet.
[void]$ws1.Range("R13").Select()

# --- Sheet "Course Handicaps": update stored selection ---
$ws2 = $wb.Worksheets.Item("Course Handicaps")
[void]$ws2.Range("T2").Select()

# --- Sheet "Course Pars": selection unchanged, just re-affirm it ---
$ws3 = $wb.Worksheets.Item("Course Pars")
[void]$ws3.Range("T2").Select()

# Leave the grid back on the first sheet, matching tabSelected state.
$ws1.Activate()
